# [#89] country in rollUp should be addressCountry
#
# The "grants" roll-up sheet is missing a "Recipient Org:Country" column.
# Insert a new column immediately before the existing
# "Recipient Org:Postal Code" column (which directly follows
# "Recipient Org:City") and label it "Recipient Org:Country". Everything
# to the right (Postal Code, Description, Web Address, Beneficiary
# Location..., Funding Org..., Grant Programme..., From an open call?,
# Related Activity, Last modified, Data Source) shifts one column right,
# and the used range grows from A1:AK1 to A1:AL1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants")

# Find the "Recipient Org:City" header so we insert right after it,
# regardless of its exact current column position.
$lastCol = $ws.Cells.Item(1, 1).End(-4161).Column
$cityCol = -1
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value() -eq "Recipient Org:City") {
        $cityCol = $c
        break
    }
}

$insertCol = $cityCol + 1
$ws.Columns.Item($insertCol).Insert()
$ws.Cells.Item(1, $insertCol).Value = "Recipient Org:Country"
